# Template_Export_Staff_Payroll.xlsx revision:
#  - add a new "Total Penghasilan" (Total Income) header in column I, pushing
#    the former I/J headers ("Total Reduction"/"Net Income") to J/K
#  - translate every header in row 1 from English to Indonesian
#  - K1 (brand new used cell) needs the same bold/centered header style as
#    the rest of row 1, so its formatting is copied over from A1
#  - re-select cell K1 (new last header) in the sheet view
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-translate / shift the header row. Columns A-H keep their column
# position; I is a brand-new header; J/K hold what used to be in I/J.
$ws.Range("A1").Value = "No"
$ws.Range("B1").Value = "Nama"
$ws.Range("C1").Value = "Tanggal Penggajian"
$ws.Range("D1").Value = "Cabang"
$ws.Range("E1").Value = "Penghasilan Pokok"
$ws.Range("F1").Value = "Insentif Kenaikan Tahunan"
$ws.Range("G1").Value = "Tidak Masuk Kerja"
$ws.Range("H1").Value = "Keterlambatan"
$ws.Range("I1").Value = "Total Penghasilan"
$ws.Range("J1").Value = "Total Pengurangan"
$ws.Range("K1").Value = "Penerimaan Bersih"

# K1 previously held no data (it was outside the used range), so it has no
# header formatting yet - copy the bold/centered style used by the rest of
# row 1 onto it.
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Approximate the post-edit column widths (author widened several columns,
# presumably to fit the longer Indonesian header text). Columns that did not
# change width (A,B,D,F,I,J,L,M) are intentionally left untouched.
$ws.Columns.Item(3).ColumnWidth = 17.58
$ws.Columns.Item(5).ColumnWidth = 20.33
$ws.Columns.Item(7).ColumnWidth = 17
$ws.Columns.Item(8).ColumnWidth = 16.08
$ws.Columns.Item(11).ColumnWidth = 19.83

# Scroll the sheet so column D is the left-most visible column, and select
# the last header cell, matching the author's final on-screen view.
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("K1").Select() | Out-Null
